$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values are written as exact text (avoid Excel numeric auto-conversion)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '70.053.53'
$ws.Range("E2").Value = '  -1.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.711.87'
$ws.Range("E3").Value = '  -1.70%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '618.96'
$ws.Range("E5").Value = '  +0.66%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.29'
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.709.26'
$ws.Range("E7").Value = '  -1.73%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.533'
$ws.Range("E9").Value = '  -1.98%  '
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.31'
$ws.Range("E11").Value = '  -2.46%  '
$ws.Range("E12").Value = '  -3.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.44'
$ws.Range("E13").Value = '  -0.54%  '
$ws.Range("E14").Value = '  -0.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.334.25'
$ws.Range("E15").Value = '  -1.74%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.713.99'
$ws.Range("E16").Value = '  -1.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '70.021.00'
$ws.Range("E17").Value = '  -1.91%  '
$ws.Range("E18").Value = '  -1.85%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.57'
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.56'
$ws.Range("E20").Value = '  -1.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '504.15'
$ws.Range("E21").Value = '  -4.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.21'
$ws.Range("E22").Value = '  -1.66%  '
$ws.Range("E23").Value = '  -4.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.59'
$ws.Range("E24").Value = '  +3.09%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.22'
$ws.Range("E25").Value = '  -2.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.48'
$ws.Range("E26").Value = '  +3.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '13.02'
$ws.Range("E27").Value = '  -4.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000133'
$ws.Range("E28").Value = '  +14.08%  '
$ws.Range("E29").Value = '  -0.29%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.48'
$ws.Range("E30").Value = '  -2.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.93'
$ws.Range("E31").Value = '  +0.59%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.88'
$ws.Range("E32").Value = '  -2.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '30.76'
$ws.Range("E33").Value = '  -4.39%  '
$ws.Range("E34").Value = '  -1.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("E36").Value = '  -0.55%  '
$ws.Range("E37").Value = '  -0.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.138'
$ws.Range("E38").Value = '  +4.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.341'
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.08'
$ws.Range("E40").Value = '  -6.80%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '430.87'
$ws.Range("E43").Value = '  -0.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.95'
$ws.Range("E44").Value = '  +5.63%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.64'
$ws.Range("E45").Value = '  -2.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.961.48'
$ws.Range("E46").Value = '  -6.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0362'
$ws.Range("E47").Value = '  -1.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.46'
$ws.Range("E48").Value = '  -2.08%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '136.21'
$ws.Range("E50").Value = '  -3.11%  '
$ws.Range("E51").Value = '  +1.48%  '

# Rows 41 and 42 swapped content (OKB <-> Arweave) with updated values
$ws.Range("B41").Value = 'Arweave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '45.87'
$ws.Range("E41").Value = '  +4.81%  '
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '49.87'
$ws.Range("E42").Value = '  -3.65%  '
